# Added category columns for CaseStatusType and ChargeDispositionType.
#
# New column C on the "CaseStatusType" and "ChargeDispositionType" lookup
# sheets carries a category grouping for each code value.

$wb = $excel.ActiveWorkbook

$wsCaseStatus  = $wb.Worksheets.Item("CaseStatusType")
$wsChargeDisp  = $wb.Worksheets.Item("ChargeDispositionType")

# --- ChargeDispositionType: fill the new category values down column C ---
# (Header C1 is filled in last, below, matching the order the unique
# shared-string table was populated in the authored workbook.)
$wsChargeDisp.Range("C2").Value  = "Charge Disposition Category A"
$wsChargeDisp.Range("C3").Value  = "Charge Disposition Category A"
$wsChargeDisp.Range("C4").Value  = "Charge Disposition Category A"
$wsChargeDisp.Range("C5").Value  = "Charge Disposition Category B"
$wsChargeDisp.Range("C6").Value  = "Charge Disposition Category B"
$wsChargeDisp.Range("C7").Value  = "Charge Disposition Category B"
$wsChargeDisp.Range("C8").Value  = "Charge Disposition Category B"
$wsChargeDisp.Range("C9").Value  = "Charge Disposition Category C"
$wsChargeDisp.Range("C10").Value = "Charge Disposition Category C"
$wsChargeDisp.Range("C11").Value = "Charge Disposition Category C"
$wsChargeDisp.Range("C12").Value = "None"
$wsChargeDisp.Range("C13").Value = "Unknown"

# --- CaseStatusType: header + category values down column C ---
$wsCaseStatus.Range("C1").Value  = "CaseStatusTypeCategory"
$wsCaseStatus.Range("C2").Value  = "Case Status Category A"
$wsCaseStatus.Range("C3").Value  = "Case Status Category A"
$wsCaseStatus.Range("C4").Value  = "Case Status Category A"
$wsCaseStatus.Range("C5").Value  = "Case Status Category A"
$wsCaseStatus.Range("C6").Value  = "Case Status Category A"
$wsCaseStatus.Range("C7").Value  = "Case Status Category A"
$wsCaseStatus.Range("C8").Value  = "Case Status Category B"
$wsCaseStatus.Range("C9").Value  = "Case Status Category B"
$wsCaseStatus.Range("C10").Value = "Case Status Category B"
$wsCaseStatus.Range("C11").Value = "Case Status Category B"
$wsCaseStatus.Range("C12").Value = "None"
$wsCaseStatus.Range("C13").Value = "Unknown"

# Re-write B13 on CaseStatusType (still "Unknown" - touched in the source
# edit alongside the new C13 cell) so the shared-string reference is kept.
$wsCaseStatus.Range("B13").Value = "Unknown"

# --- ChargeDispositionType: header last ---
$wsChargeDisp.Range("C1").Value = "ChargeDispositionTypeCategory"

# --- Column widths for the new column C on both sheets ---
$wsCaseStatus.Columns.Item(3).ColumnWidth = 30.8
$wsChargeDisp.Columns.Item(3).ColumnWidth = 41.5

# --- Selections / active sheet to match the authored workbook state ---
$wsCaseStatus.Range("C16").Select() | Out-Null
$wsChargeDisp.Activate() | Out-Null
$wsChargeDisp.Range("C2").Select() | Out-Null
